# Weekly refresh: a new price record is added at the top of the data
# block (row 69), pushing the existing rows 69-82 down to 70-83. The
# record that used to be the last one (row 82) becomes the new last
# row (83).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 69; this shifts rows 69:82 down to
# 70:83 (carrying their values/styles with them) and grows the used
# range / dimension to A1:R83 automatically.
$ws.Rows.Item(69).Insert()

# Populate the newly-inserted row 69 with this week's record.
$ws.Cells.Item(69, 1).Value = 7
$ws.Cells.Item(69, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(69, 3).Value = "Ñuble"
$ws.Cells.Item(69, 4).Value = 45077
$ws.Cells.Item(69, 5).Value = 16
$ws.Cells.Item(69, 6).Value = 100112013
$ws.Cells.Item(69, 7).Value = "Alcachofa"
$ws.Cells.Item(69, 8).Value = "Española"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 60
$ws.Cells.Item(69, 11).Value = 15000
$ws.Cells.Item(69, 12).Value = 16000
$ws.Cells.Item(69, 13).Value = 15500
$ws.Cells.Item(69, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(69, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(69, 16).Value = 517
$ws.Cells.Item(69, 17).Value = 30
$ws.Cells.Item(69, 18).Value = "Hortaliza"
